$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 546.7765096666667
$ws.Range("H2").Value = 1640.329529
$ws.Range("I2").Value = 0.6285526459909564
$ws.Range("J2").Value = 0.6285526459909564
$ws.Range("M2").Value = 71.541692
$ws.Range("N2").Value = 214.625076
$ws.Range("O2").Value = 0.6133071420247926
$ws.Range("P2").Value = 0.6133071420247926
$ws.Range("Q2").Value = 39117.31664740769
$ws.Range("R2").Value = 352055.8498266692
$ws.Range("S2").Value = 0.3854958269248347
$ws.Range("T2").Value = 0.3854958269248347
$ws.Range("G3").Value = 546.7765096666667
$ws.Range("H3").Value = 1640.329529
$ws.Range("I3").Value = 0.6285526459909564
$ws.Range("J3").Value = 0.6285526459909564
$ws.Range("O3").Value = 0.08457024278578675
$ws.Range("P3").Value = 0.08457024278578675
$ws.Range("Q3").Value = 5393.971045369034
$ws.Range("R3").Value = 48545.73940832131
$ws.Range("S3").Value = 0.05315684987510385
$ws.Range("T3").Value = 0.05315684987510385
$ws.Range("G4").Value = 546.7765096666667
$ws.Range("H4").Value = 1640.329529
$ws.Range("I4").Value = 0.6285526459909564
$ws.Range("J4").Value = 0.6285526459909564
$ws.Range("M4").Value = 35.05835333333334
$ws.Range("N4").Value = 105.17506
$ws.Range("O4").Value = 0.3005455684073286
$ws.Range("P4").Value = 0.3005455684073286
$ws.Range("Q4").Value = 19169.08407026075
$ws.Range("R4").Value = 172521.7566323467
$ws.Range("S4").Value = 0.1889087122632824
$ws.Range("T4").Value = 0.1889087122632824
$ws.Range("G5").Value = 546.7765096666667
$ws.Range("H5").Value = 1640.329529
$ws.Range("I5").Value = 0.6285526459909564
$ws.Range("J5").Value = 0.6285526459909564
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.183961
$ws.Range("N5").Value = 0.551883
$ws.Range("O5").Value = 0.001577046782092083
$ws.Range("P5").Value = 0.001577046782092083
$ws.Range("Q5").Value = 100.5855534947897
$ws.Range("R5").Value = 905.2699814531071
$ws.Range("S5").Value = 0.0009912569277355017
$ws.Range("T5").Value = 0.0009912569277355017
$ws.Range("I6").Value = 0.1861770314550556
$ws.Range("J6").Value = 0.1861770314550556
$ws.Range("M6").Value = 71.541692
$ws.Range("N6").Value = 214.625076
$ws.Range("O6").Value = 0.6133071420247926
$ws.Range("P6").Value = 0.6133071420247926
$ws.Range("Q6").Value = 11586.53286777601
$ws.Range("R6").Value = 104278.7958099841
$ws.Range("S6").Value = 0.1141837030723601
$ws.Range("T6").Value = 0.1141837030723601
$ws.Range("I7").Value = 0.1861770314550556
$ws.Range("J7").Value = 0.1861770314550556
$ws.Range("O7").Value = 0.08457024278578675
$ws.Range("P7").Value = 0.08457024278578675
$ws.Range("S7").Value = 0.01574503675129111
$ws.Range("T7").Value = 0.01574503675129111
$ws.Range("I8").Value = 0.1861770314550556
$ws.Range("J8").Value = 0.1861770314550556
$ws.Range("M8").Value = 35.05835333333334
$ws.Range("N8").Value = 105.17506
$ws.Range("O8").Value = 0.3005455684073286
$ws.Range("P8").Value = 0.3005455684073286
$ws.Range("Q8").Value = 5677.874702582813
$ws.Range("R8").Value = 51100.87232324532
$ws.Range("S8").Value = 0.05595468174304879
$ws.Range("T8").Value = 0.05595468174304879
$ws.Range("I9").Value = 0.1861770314550556
$ws.Range("J9").Value = 0.1861770314550556
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.183961
$ws.Range("N9").Value = 0.551883
$ws.Range("O9").Value = 0.001577046782092083
$ws.Range("P9").Value = 0.001577046782092083
$ws.Range("Q9").Value = 29.793398972014
$ws.Range("R9").Value = 268.140590748126
$ws.Range("S9").Value = 0.0002936098883556519
$ws.Range("T9").Value = 0.0002936098883556519
$ws.Range("G10").Value = 160.630483
$ws.Range("H10").Value = 481.891449
$ws.Range("I10").Value = 0.1846544489960017
$ws.Range("J10").Value = 0.1846544489960017
$ws.Range("M10").Value = 71.541692
$ws.Range("N10").Value = 214.625076
$ws.Range("O10").Value = 0.6133071420247926
$ws.Range("P10").Value = 0.6133071420247926
$ws.Range("Q10").Value = 11491.77654059723
$ws.Range("R10").Value = 103425.9888653751
$ws.Range("S10").Value = 0.1132498923759006
$ws.Range("T10").Value = 0.1132498923759006
$ws.Range("G11").Value = 160.630483
$ws.Range("H11").Value = 481.891449
$ws.Range("I11").Value = 0.1846544489960017
$ws.Range("J11").Value = 0.1846544489960017
$ws.Range("O11").Value = 0.08457024278578675
$ws.Range("P11").Value = 0.08457024278578675
$ws.Range("Q11").Value = 1584.625818753354
$ws.Range("R11").Value = 14261.63236878019
$ws.Range("S11").Value = 0.01561627158306754
$ws.Range("T11").Value = 0.01561627158306754
$ws.Range("G12").Value = 160.630483
$ws.Range("H12").Value = 481.891449
$ws.Range("I12").Value = 0.1846544489960017
$ws.Range("J12").Value = 0.1846544489960017
$ws.Range("M12").Value = 35.05835333333334
$ws.Range("N12").Value = 105.17506
$ws.Range("O12").Value = 0.3005455684073286
$ws.Range("P12").Value = 0.3005455684073286
$ws.Range("Q12").Value = 5631.440229117994
$ws.Range("R12").Value = 50682.96206206194
$ws.Range("S12").Value = 0.0554970763324454
$ws.Range("T12").Value = 0.0554970763324454
$ws.Range("G13").Value = 160.630483
$ws.Range("H13").Value = 481.891449
$ws.Range("I13").Value = 0.1846544489960017
$ws.Range("J13").Value = 0.1846544489960017
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.183961
$ws.Range("N13").Value = 0.551883
$ws.Range("O13").Value = 0.001577046782092083
$ws.Range("P13").Value = 0.001577046782092083
$ws.Range("Q13").Value = 29.549744283163
$ws.Range("R13").Value = 265.947698548467
$ws.Range("S13").Value = 0.0002912087045881311
$ws.Range("T13").Value = 0.0002912087045881311
$ws.Range("G14").Value = 0.5357470000000001
$ws.Range("H14").Value = 1.607241
$ws.Range("I14").Value = 0.0006158735579862568
$ws.Range("J14").Value = 0.0006158735579862568
$ws.Range("M14").Value = 71.541692
$ws.Range("N14").Value = 214.625076
$ws.Range("O14").Value = 0.6133071420247926
$ws.Range("P14").Value = 0.6133071420247926
$ws.Range("Q14").Value = 38.32824686392401
$ws.Range("R14").Value = 344.954221775316
$ws.Range("S14").Value = 0.0003777196516971915
$ws.Range("T14").Value = 0.0003777196516971915
$ws.Range("G15").Value = 0.5357470000000001
$ws.Range("H15").Value = 1.607241
$ws.Range("I15").Value = 0.0006158735579862568
$ws.Range("J15").Value = 0.0006158735579862568
$ws.Range("O15").Value = 0.08457024278578675
$ws.Range("P15").Value = 0.08457024278578675
$ws.Range("Q15").Value = 5.285164513386001
$ws.Range("R15").Value = 47.566480620474
$ws.Range("S15").Value = 0.00005208457632424404
$ws.Range("T15").Value = 0.00005208457632424404
$ws.Range("G16").Value = 0.5357470000000001
$ws.Range("H16").Value = 1.607241
$ws.Range("I16").Value = 0.0006158735579862568
$ws.Range("J16").Value = 0.0006158735579862568
$ws.Range("M16").Value = 35.05835333333334
$ws.Range("N16").Value = 105.17506
$ws.Range("O16").Value = 0.3005455684073286
$ws.Range("P16").Value = 0.3005455684073286
$ws.Range("Q16").Value = 18.78240762327334
$ws.Range("R16").Value = 169.04166860946
$ws.Range("S16").Value = 0.0001850980685520234
$ws.Range("T16").Value = 0.0001850980685520234
$ws.Range("G17").Value = 0.5357470000000001
$ws.Range("H17").Value = 1.607241
$ws.Range("I17").Value = 0.0006158735579862568
$ws.Range("J17").Value = 0.0006158735579862568
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.183961
$ws.Range("N17").Value = 0.551883
$ws.Range("O17").Value = 0.001577046782092083
$ws.Range("P17").Value = 0.001577046782092083
$ws.Range("Q17").Value = 0.09855655386700002
$ws.Range("R17").Value = 0.8870089848030001
$ws.Range("S17").Value = 0.0000009712614127978279
$ws.Range("T17").Value = 0.0000009712614127978279
